# Auto-generated update of Leve profit calculations across multiple sheets.
# Values mirror a refreshed Universalis market-price pull for the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2247.2727
$ws.Range("I6").Value = 246.66667
$ws.Range("J6").Value = 2997.5
$ws.Range("K6").Value = 740.00001
$ws.Range("L6").Value = 8992.5
$ws.Range("M6").Value = -628.00001
$ws.Range("N6").Value = -9216.5

$ws.Range("H64").Value = 3490
$ws.Range("I64").Value = 3326.6667
$ws.Range("J64").Value = 3653.3333
$ws.Range("K64").Value = 3326.6667
$ws.Range("L64").Value = 3653.3333
$ws.Range("M64").Value = -3078.6667
$ws.Range("N64").Value = -4149.3333

$ws.Range("H67").Value = 3490
$ws.Range("I67").Value = 3326.6667
$ws.Range("J67").Value = 3653.3333
$ws.Range("K67").Value = 3326.6667
$ws.Range("L67").Value = 3653.3333
$ws.Range("M67").Value = -2468.6667
$ws.Range("N67").Value = -5369.3333

$ws.Range("H74").Value = 4416.625
$ws.Range("I74").Value = 3896.5
$ws.Range("J74").Value = 4590
$ws.Range("K74").Value = 3896.5
$ws.Range("L74").Value = 4590
$ws.Range("M74").Value = -2960.5
$ws.Range("N74").Value = -6462

$ws.Range("H76").Value = 2896.8235
$ws.Range("I76").Value = 2831.5
$ws.Range("K76").Value = 2831.5
$ws.Range("M76").Value = -2516.5

$ws.Range("H77").Value = 4416.625
$ws.Range("I77").Value = 3896.5
$ws.Range("J77").Value = 4590
$ws.Range("K77").Value = 19482.5
$ws.Range("L77").Value = 22950
$ws.Range("M77").Value = -14802.5
$ws.Range("N77").Value = -32310

$ws.Range("H79").Value = 2896.8235
$ws.Range("I79").Value = 2831.5
$ws.Range("K79").Value = 2831.5
$ws.Range("M79").Value = -1739.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6929.78
$ws.Range("I32").Value = 6472
$ws.Range("J32").Value = 9993.385
$ws.Range("K32").Value = 6472
$ws.Range("L32").Value = 9993.385
$ws.Range("M32").Value = -6185
$ws.Range("N32").Value = -10567.385

$ws.Range("H132").Value = 1857.7727
$ws.Range("I132").Value = 1326.7693
$ws.Range("K132").Value = 3980.3079
$ws.Range("M132").Value = -1450.3079

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2819.84
$ws.Range("I31").Value = 1659.359
$ws.Range("J31").Value = 3561.7869
$ws.Range("K31").Value = 1659.359
$ws.Range("L31").Value = 3561.7869
$ws.Range("M31").Value = -1364.359
$ws.Range("N31").Value = -4151.7869

$ws.Range("H34").Value = 2819.84
$ws.Range("I34").Value = 1659.359
$ws.Range("J34").Value = 3561.7869
$ws.Range("K34").Value = 1659.359
$ws.Range("L34").Value = 3561.7869
$ws.Range("M34").Value = -1457.359
$ws.Range("N34").Value = -3965.7869

$ws.Range("H131").Value = 16812.666
$ws.Range("J131").Value = 17277.857
$ws.Range("L131").Value = 17277.857
$ws.Range("N131").Value = -27357.857

$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 35.933334
$ws.Range("I2").Value = 23.222221
$ws.Range("J2").Value = 55
$ws.Range("K2").Value = 139.333326
$ws.Range("L2").Value = 330
$ws.Range("M2").Value = -26.333326
$ws.Range("N2").Value = -556

$ws.Range("H7").Value = 393
$ws.Range("I7").Value = 90.5
$ws.Range("J7").Value = 468.625
$ws.Range("K7").Value = 271.5
$ws.Range("L7").Value = 1405.875
$ws.Range("M7").Value = -159.5
$ws.Range("N7").Value = -1629.875

$ws.Range("H15").Value = 232.5
$ws.Range("I15").Value = 143.33333
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 429.99999
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = -289.99999
$ws.Range("N15").Value = -1780

$ws.Range("H34").Value = 10941.6
$ws.Range("I34").Value = 188
$ws.Range("J34").Value = 21695.2
$ws.Range("K34").Value = 564
$ws.Range("L34").Value = 65085.60000000001
$ws.Range("M34").Value = -480
$ws.Range("N34").Value = -65253.60000000001

$ws.Range("H62").Value = 16942.857
$ws.Range("I62").Value = 600
$ws.Range("J62").Value = 19666.666
$ws.Range("K62").Value = 1800
$ws.Range("L62").Value = 58999.99800000001
$ws.Range("M62").Value = -1114
$ws.Range("N62").Value = -60371.99800000001

$ws.Range("H65").Value = 16942.857
$ws.Range("I65").Value = 600
$ws.Range("J65").Value = 19666.666
$ws.Range("K65").Value = 5400
$ws.Range("L65").Value = 176999.994
$ws.Range("M65").Value = -1968
$ws.Range("N65").Value = -183863.994

$ws.Range("H68").Value = 1320
$ws.Range("I68").Value = 413.33334
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1240.00002
$ws.Range("L68").Value = 6000
$ws.Range("M68").Value = -429.0000199999999
$ws.Range("N68").Value = -7622

$ws.Range("H71").Value = 1320
$ws.Range("I71").Value = 413.33334
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 3720.00006
$ws.Range("L71").Value = 18000
$ws.Range("M71").Value = 335.9999399999997
$ws.Range("N71").Value = -26112

$ws.Range("H110").Value = 3110
$ws.Range("I110").Value = 440
$ws.Range("K110").Value = 1320
$ws.Range("M110").Value = 2770

$ws.Range("H114").Value = 653.28
$ws.Range("I114").Value = 203.58333
$ws.Range("K114").Value = 610.74999
$ws.Range("M114").Value = 2643.25001

$ws.Range("H121").Value = 10143.4375
$ws.Range("I121").Value = 225.44444
$ws.Range("J121").Value = 22895.143
$ws.Range("K121").Value = 676.33332
$ws.Range("L121").Value = 68685.429
$ws.Range("M121").Value = 633.66668
$ws.Range("N121").Value = -71305.429

$ws.Range("H131").Value = 1182.1632
$ws.Range("I131").Value = 2354.5454
$ws.Range("J131").Value = 1033.931
$ws.Range("K131").Value = 7063.6362
$ws.Range("L131").Value = 3101.793
$ws.Range("M131").Value = -2023.6362
$ws.Range("N131").Value = -13181.793

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4887
$ws.Range("I70").Value = 4675.7
$ws.Range("J70").Value = 7000
$ws.Range("K70").Value = 4675.7
$ws.Range("L70").Value = 7000
$ws.Range("M70").Value = -4405.7
$ws.Range("N70").Value = -7540

$ws.Range("H73").Value = 4887
$ws.Range("I73").Value = 4675.7
$ws.Range("J73").Value = 7000
$ws.Range("K73").Value = 4675.7
$ws.Range("L73").Value = 7000
$ws.Range("M73").Value = -3739.7
$ws.Range("N73").Value = -8872

$ws.Range("H126").Value = 591783.1
$ws.Range("I126").Value = 2255.4443
$ws.Range("J126").Value = 1255001.8
$ws.Range("K126").Value = 6766.3329
$ws.Range("L126").Value = 3765005.4
$ws.Range("M126").Value = -4296.3329
$ws.Range("N126").Value = -3769945.4

$ws.Range("H132").Value = 3213.6511
$ws.Range("I132").Value = 3177.2693
$ws.Range("K132").Value = 9531.8079
$ws.Range("M132").Value = -7001.8079

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2303.2942
$ws.Range("I132").Value = 1631.8
$ws.Range("J132").Value = 3772.1875
$ws.Range("K132").Value = 4895.4
$ws.Range("L132").Value = 11316.5625
$ws.Range("M132").Value = -2365.4
$ws.Range("N132").Value = -16376.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 15620.281
$ws.Range("I132").Value = 2819.1738
$ws.Range("J132").Value = 48334.223
$ws.Range("K132").Value = 8457.5214
$ws.Range("L132").Value = 145002.669
$ws.Range("M132").Value = -5927.5214
$ws.Range("N132").Value = -150062.669
